$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels for columns J and K
$ws.Range("J1").Value = "進貨價"
$ws.Range("K1").Value = "url"

# New J (進貨價) values per row, taken from target diff
$jValues = @{
    2  = 45.15
    3  = 94.89
    4  = 29.73
    5  = 32.44
    6  = 83.42
    7  = 96.84
    8  = 81.73999999999999
    9  = 39.98
    10 = 24.66
    11 = 136.3
    12 = 63.56
    13 = 30.96
    14 = 12.05
    15 = 72.81999999999999
    16 = 13.88
    17 = 119.04
    18 = 7.54
    19 = 13.18
    20 = 49.43
    21 = 67.08
}

for ($row = 2; $row -le 21; $row++) {
    $url = $ws.Range("J$row").Value2
    $ws.Range("K$row").Value = $url
    $ws.Range("J$row").Value = $jValues[$row]
}
